$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every data row.
# It is stamped with the current date whenever the report is regenerated,
# so every existing row moves from 45181 (2023-09-12) to 45182 (2023-09-13).
# The sheet has data in rows 2..451 (row 1 is the header).
for ($r = 2; $r -le 451; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = 45182
    }
}
